# Append the new data row (row 65) to Sheet1, as produced by the
# 2026-01-28 run of the profit-tracking job.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 65

# Column A holds the date as literal text (matching the existing rows),
# so use a leading apostrophe to stop Excel's automatic date conversion,
# then reset the cell style back to Normal (quote-prefix doesn't change
# the visible formatting, only how the value was entered).
$ws.Cells.Item($row, 1).Value  = "'01/28/2026"
$ws.Cells.Item($row, 1).Style  = "Normal"
$ws.Cells.Item($row, 2).Value  = 11625.89
$ws.Cells.Item($row, 3).Value  = 0.2434540290739574
$ws.Cells.Item($row, 4).Value  = 0.7565459709260426
$ws.Cells.Item($row, 5).Value  = -205.6
$ws.Cells.Item($row, 6).Value  = -27.97
$ws.Cells.Item($row, 7).Value  = -21996.24
$ws.Cells.Item($row, 8).Value  = -71.44
$ws.Cells.Item($row, 9).Value  = -422.48
$ws.Cells.Item($row, 10).Value = -12.99
$ws.Cells.Item($row, 11).Value = -22418.72
$ws.Cells.Item($row, 12).Value = -65.84999999999999
